$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1491.6666
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1491.6666
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 4474.9998
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -4774.9998

$ws.Range("H105").Value = 14055
$ws.Range("J105").Value = 15203.333
$ws.Range("L105").Value = 15203.333
$ws.Range("N105").Value = -22191.333

$ws.Range("H132").Value = 1561.7
$ws.Range("I132").Value = 1152.125
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 3456.375
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -926.375
$ws.Range("N132").Value = -14660

$ws.Range("H137").Value = 4115.1665
$ws.Range("I137").Value = 4115.1665
$ws.Range("K137").Value = 12345.4995
$ws.Range("M137").Value = -9795.499500000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 736.4
$ws.Range("I32").Value = 736.4
$ws.Range("K32").Value = 736.4
$ws.Range("M32").Value = -449.4

$ws.Range("H74").Value = 3210.05
$ws.Range("I74").Value = 3221.1052
$ws.Range("K74").Value = 3221.1052
$ws.Range("M74").Value = -2347.1052

$ws.Range("H77").Value = 3210.05
$ws.Range("I77").Value = 3221.1052
$ws.Range("K77").Value = 16105.526
$ws.Range("M77").Value = -11737.526

$ws.Range("H86").Value = 45157
$ws.Range("I86").Value = 30000
$ws.Range("J86").Value = 60314
$ws.Range("K86").Value = 30000
$ws.Range("L86").Value = 60314
$ws.Range("M86").Value = -28814
$ws.Range("N86").Value = -62686

$ws.Range("H89").Value = 45157
$ws.Range("I89").Value = 30000
$ws.Range("J89").Value = 60314
$ws.Range("K89").Value = 90000
$ws.Range("L89").Value = 180942
$ws.Range("M89").Value = -84072
$ws.Range("N89").Value = -192798

$ws.Range("H92").Value = 49999.5
$ws.Range("J92").Value = 49999.5
$ws.Range("L92").Value = 49999.5
$ws.Range("N92").Value = -54991.5

$ws.Range("H101").Value = 224325.25
$ws.Range("J101").Value = 224325.25
$ws.Range("L101").Value = 224325.25
$ws.Range("N101").Value = -230815.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4002.1428
$ws.Range("I20").Value = 3603
$ws.Range("K20").Value = 3603
$ws.Range("M20").Value = -3356

$ws.Range("H22").Value = 550
$ws.Range("I22").Value = 550
$ws.Range("K22").Value = 550
$ws.Range("M22").Value = -377

$ws.Range("H86").Value = 3666.3333
$ws.Range("I86").Value = 2399.6
$ws.Range("K86").Value = 2399.6
$ws.Range("M86").Value = -1276.6

$ws.Range("H89").Value = 3666.3333
$ws.Range("I89").Value = 2399.6
$ws.Range("K89").Value = 11998
$ws.Range("M89").Value = -6382

$ws.Range("H92").Value = 91799.8
$ws.Range("J92").Value = 91799.8
$ws.Range("L92").Value = 91799.8
$ws.Range("N92").Value = -96791.8

$ws.Range("H134").Value = 6953
$ws.Range("I134").Value = 6039.4
$ws.Range("K134").Value = 18118.2
$ws.Range("M134").Value = -15583.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 75.888885
$ws.Range("I7").Value = 18.066668
$ws.Range("J7").Value = 365
$ws.Range("K7").Value = 18.066668
$ws.Range("L7").Value = 365
$ws.Range("M7").Value = 94.93333200000001
$ws.Range("N7").Value = -591

$ws.Range("H8").Value = 2499
$ws.Range("J8").Value = 4989
$ws.Range("L8").Value = 4989
$ws.Range("N8").Value = -5269

$ws.Range("H16").Value = 924.6667
$ws.Range("I16").Value = 924.6667
$ws.Range("K16").Value = 924.6667
$ws.Range("M16").Value = -637.6667

$ws.Range("H58").Value = 1242.5
$ws.Range("I58").Value = 485
$ws.Range("K58").Value = 485
$ws.Range("M58").Value = -282

$ws.Range("H113").Value = 924.6667
$ws.Range("I113").Value = 924.6667
$ws.Range("K113").Value = 924.6667
$ws.Range("M113").Value = 1245.3333

$ws.Range("H136").Value = 1242.5
$ws.Range("I136").Value = 485
$ws.Range("K136").Value = 1455
$ws.Range("M136").Value = 1095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 364.35715
$ws.Range("I6").Value = 26.166666
$ws.Range("J6").Value = 456.5909
$ws.Range("K6").Value = 78.49999800000001
$ws.Range("L6").Value = 1369.7727
$ws.Range("M6").Value = 34.50000199999999
$ws.Range("N6").Value = -1595.7727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 19113.666
$ws.Range("J104").Value = 19113.666
$ws.Range("L104").Value = 19113.666
$ws.Range("N104").Value = -26101.666

$ws.Range("H107").Value = 634.8
$ws.Range("I107").Value = 700
$ws.Range("J107").Value = 374
$ws.Range("K107").Value = 700
$ws.Range("L107").Value = 374
$ws.Range("M107").Value = 1220
$ws.Range("N107").Value = -4214

$ws.Range("H132").Value = 1349
$ws.Range("I132").Value = 1277.25
$ws.Range("J132").Value = 1492.5
$ws.Range("K132").Value = 3831.75
$ws.Range("L132").Value = 4477.5
$ws.Range("M132").Value = -1301.75
$ws.Range("N132").Value = -9537.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1716.5
$ws.Range("I22").Value = 1249.5
$ws.Range("K22").Value = 1249.5
$ws.Range("M22").Value = -954.5

$ws.Range("H27").Value = 1716.5
$ws.Range("I27").Value = 1249.5
$ws.Range("K27").Value = 1249.5
$ws.Range("M27").Value = -1142.5

$ws.Range("H44").Value = 12500
$ws.Range("J44").Value = 12500
$ws.Range("L44").Value = 12500
$ws.Range("N44").Value = -13412

$ws.Range("H46").Value = 334946.34
$ws.Range("I46").Value = 667393.3
$ws.Range("K46").Value = 667393.3
$ws.Range("M46").Value = -667205.3

$ws.Range("H61").Value = 100005
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H113").Value = 100005
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H132").Value = 999.3333
$ws.Range("I132").Value = 999.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2997.9999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -467.9998999999998
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1699.5
$ws.Range("I136").Value = 1699.5
$ws.Range("K136").Value = 5098.5
$ws.Range("M136").Value = -2548.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 4040642.2
$ws.Range("I10").Value = 6704399.5
$ws.Range("J10").Value = 45006
$ws.Range("K10").Value = 6704399.5
$ws.Range("L10").Value = 45006
$ws.Range("M10").Value = -6704230.5
$ws.Range("N10").Value = -45344

$ws.Range("H117").Value = 39998
$ws.Range("J117").Value = 39998
$ws.Range("L117").Value = 39998
$ws.Range("N117").Value = -49176

$ws.Range("H132").Value = 2680
$ws.Range("I132").Value = 2680
$ws.Range("K132").Value = 8040
$ws.Range("M132").Value = -5510

$ws.Range("H136").Value = 10491.929
$ws.Range("I136").Value = 10265.333
$ws.Range("K136").Value = 30795.999
$ws.Range("M136").Value = -28245.999
